$d = $word.ActiveDocument

$replacements = @(
    @{old="397÷4="; new="731÷5="},
    @{old="911÷2="; new="301÷7="},
    @{old="150÷9="; new="487÷2="},
    @{old="826÷2="; new="778÷7="},
    @{old="807÷9="; new="948÷6="},
    @{old="491÷5="; new="289÷4="},
    @{old="744÷9="; new="490÷8="},
    @{old="130÷2="; new="564÷4="},
    @{old="782÷2="; new="550÷7="},
    @{old="539÷9="; new="164÷9="},
    @{old="552÷2="; new="872÷2="},
    @{old="213÷5="; new="151÷9="},
    @{old="270÷2="; new="165÷6="},
    @{old="400÷3="; new="225÷3="},
    @{old="913÷6="; new="508÷5="},
    @{old="415÷2="; new="538÷3="},
    @{old="821÷9="; new="668÷4="},
    @{old="963÷4="; new="507÷6="},
    @{old="426÷9="; new="674÷5="},
    @{old="763÷5="; new="222÷7="},
    @{old="246÷5="; new="491÷4="},
    @{old="661÷9="; new="330÷4="},
    @{old="796÷4="; new="722÷9="},
    @{old="991÷2="; new="157÷9="},
    @{old="928÷3="; new="992÷4="}
)

foreach ($r in $replacements) {
    $rng = $d.Content
    $rng.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
